$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.678.14"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "2.465.85"
$ws.Range("E3").Value = "  -1.10%  "
$ws.Range("E4").Value = "  +0.53%  "
$ws.Range("D5").Value = "'315.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.66%  "
$ws.Range("D6").Value = "'92.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.51%  "
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("E8").Value = "  +0.44%  "
$ws.Range("D9").Value = "'0.513"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'32.53"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.14%  "
$ws.Range("D11").Value = "'0.0836"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.35%  "
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").Value = "2.844.41"
$ws.Range("E13").Value = "  -1.08%  "
$ws.Range("D14").Value = "'6.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").Value = "'15.81"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.71%  "
$ws.Range("D16").Value = "2.460.54"
$ws.Range("E16").Value = "  -2.44%  "
$ws.Range("D17").Value = "'0.779"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.32%  "
$ws.Range("D18").Value = "41.655.20"
$ws.Range("E18").Value = "  -0.02%  "
$ws.Range("E19").Value = "  +2.25%  "
$ws.Range("D20").Value = "0.0₃0944"
$ws.Range("E20").Value = "  +2.15%  "
$ws.Range("D21").Value = "'70.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("D22").Value = "'11.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.93%  "
$ws.Range("D23").Value = "'238.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.95%  "
$ws.Range("E24").Value = "  -0.30%  "
$ws.Range("E25").Value = "  +0.40%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").Value = "'24.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.16%  "
$ws.Range("E28").Value = "  +0.51%  "
$ws.Range("E29").Value = "  +0.53%  "
$ws.Range("D30").Value = "'35.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.76%  "
$ws.Range("D31").Value = "'155.64"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.67%  "
$ws.Range("D32").Value = "'5.51"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.66%  "
$ws.Range("E33").Value = "  +0.37%  "
$ws.Range("D34").Value = "'0.0761"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.50%  "
$ws.Range("E35").Value = "  -0.43%  "
$ws.Range("D36").Value = "'17.45"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.33%  "
$ws.Range("E37").Value = "  -2.36%  "
$ws.Range("E38").Value = "  +0.97%  "
$ws.Range("E39").Value = "  +0.32%  "
$ws.Range("D40").Value = "'1.80"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.31%  "
$ws.Range("E41").Value = "  -5.34%  "
$ws.Range("E42").Value = "  +0.51%  "
$ws.Range("D43").Value = "1.977.05"
$ws.Range("E43").Value = "  +1.32%  "
$ws.Range("D44").Value = "'18.86"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.66%  "
$ws.Range("E45").Value = "  -1.04%  "
$ws.Range("E46").Value = "  -1.59%  "
$ws.Range("D47").Value = "'9.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.97%  "
$ws.Range("D48").Value = "2.701.32"
$ws.Range("E48").Value = "  -1.25%  "
$ws.Range("D49").Value = "'96.81"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("D50").Value = "'67.04"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.03%  "
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").Value = "'72.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.44%  "
